$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) cells whose new value is unambiguously text (contains two
# "." separators, e.g. thousands + decimal), so a plain .Value assignment keeps it text.
$priceTextValues = @{
    "D2" = "27.407.25"
    "D3" = "1.641.49"
    "D12" = "1.873.57"
    "D13" = "1.640.26"
    "D17" = "27.384.95"
    "D34" = "1.399.09"
    "D46" = "1.784.25"
}
foreach ($addr in $priceTextValues.Keys) {
    $ws.Range($addr).Value = $priceTextValues[$addr]
}

# --- Price (column D) cells whose new value looks like a plain number (e.g. "212.00",
# "0.536"). Excel would silently convert a plain .Value assignment into a numeric value
# (dropping trailing zeros / changing representation), so force text entry by temporarily
# switching the cell to a text number format, then restore the original (default) style.
$priceNumericLookingValues = @{
    "D5" = "212.00"
    "D6" = "0.536"
    "D8" = "23.25"
    "D14" = "4.03"
    "D16" = "64.20"
    "D18" = "228.23"
    "D20" = "7.47"
    "D22" = "4.30"
    "D23" = "9.29"
    "D25" = "147.62"
    "D29" = "15.52"
    "D35" = "1.56"
    "D37" = "0.561"
    "D38" = "0.877"
    "D41" = "0.999"
    "D44" = "0.787"
    "D45" = "64.25"
    "D48" = "87.30"
    "D50" = "0.0985"
}
$defaultStyle = $ws.Range("C2").Style
foreach ($addr in $priceNumericLookingValues.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceNumericLookingValues[$addr]
    $cell.Style = $defaultStyle
}

# --- Volume(1h) (column E) cells: percentage text padded with two leading/trailing
# spaces. These never parse as numbers (because of the spaces), so plain assignment is safe.
$volumeValues = @{
    "E2" = "  -0.66%  "
    "E3" = "  -1.54%  "
    "E4" = "  -0.02%  "
    "E6" = "  +4.34%  "
    "E7" = "  +0.01%  "
    "E8" = "  -1.46%  "
    "E9" = "  -2.35%  "
    "E10" = "  -2.06%  "
    "E11" = "  +1.27%  "
    "E12" = "  -1.52%  "
    "E13" = "  -1.39%  "
    "E14" = "  -3.19%  "
    "E15" = "  +0.65%  "
    "E16" = "  -3.19%  "
    "E17" = "  -0.82%  "
    "E18" = "  -9.12%  "
    "E19" = "  -1.87%  "
    "E20" = "  -1.11%  "
    "E21" = "  +0.05%  "
    "E22" = "  -4.43%  "
    "E23" = "  -0.16%  "
    "E24" = "  +0.28%  "
    "E25" = "  +0.73%  "
    "E26" = "  +2.36%  "
    "E27" = "  -3.07%  "
    "E28" = "  -0.05%  "
    "E29" = "  -6.24%  "
    "E30" = "  -4.94%  "
    "E31" = "  -4.24%  "
    "E32" = "  -2.76%  "
    "E33" = "  -0.82%  "
    "E34" = "  -4.90%  "
    "E35" = "  -1.20%  "
    "E36" = "  -0.33%  "
    "E37" = "  -2.75%  "
    "E38" = "  -7.03%  "
    "E39" = "  -3.24%  "
    "E40" = "  -0.39%  "
    "E41" = "  +0.00%  "
    "E42" = "  +0.74%  "
    "E43" = "  +0.38%  "
    "E44" = "  -0.40%  "
    "E45" = "  -7.84%  "
    "E47" = "  -3.92%  "
    "E49" = "  -3.65%  "
    "E50" = "  -3.48%  "
}
foreach ($addr in $volumeValues.Keys) {
    $ws.Range($addr).Value = $volumeValues[$addr]
}
